$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 945-1013: data shifts down by one weekly block (3 rows) ---
# Row 945
$ws.Cells.Item(945, 4).Value = 44783
$ws.Cells.Item(945, 11).Value = 2000
$ws.Cells.Item(945, 12).Value = 2500
$ws.Cells.Item(945, 13).Value = 2250
$ws.Cells.Item(945, 16).Value = 225

# Row 946
$ws.Cells.Item(946, 4).Value = 44783
$ws.Cells.Item(946, 10).Value = 270
$ws.Cells.Item(946, 11).Value = 1500
$ws.Cells.Item(946, 12).Value = 2000
$ws.Cells.Item(946, 13).Value = 1750
$ws.Cells.Item(946, 16).Value = 175

# Row 947
$ws.Cells.Item(947, 4).Value = 44783
$ws.Cells.Item(947, 11).Value = 1000
$ws.Cells.Item(947, 12).Value = 1500
$ws.Cells.Item(947, 13).Value = 1250
$ws.Cells.Item(947, 16).Value = 125

# Row 948
$ws.Cells.Item(948, 4).Value = 44425
$ws.Cells.Item(948, 10).Value = 250
$ws.Cells.Item(948, 11).Value = 2500
$ws.Cells.Item(948, 12).Value = 3000
$ws.Cells.Item(948, 13).Value = 2750
$ws.Cells.Item(948, 16).Value = 275

# Row 949
$ws.Cells.Item(949, 4).Value = 44425
$ws.Cells.Item(949, 10).Value = 300
$ws.Cells.Item(949, 11).Value = 2300
$ws.Cells.Item(949, 12).Value = 2500
$ws.Cells.Item(949, 13).Value = 2400
$ws.Cells.Item(949, 16).Value = 240

# Row 950
$ws.Cells.Item(950, 4).Value = 44425
$ws.Cells.Item(950, 10).Value = 300
$ws.Cells.Item(950, 12).Value = 2300
$ws.Cells.Item(950, 13).Value = 2150
$ws.Cells.Item(950, 16).Value = 215

# Row 951
$ws.Cells.Item(951, 4).Value = 44377
$ws.Cells.Item(951, 10).Value = 300
$ws.Cells.Item(951, 11).Value = 3000
$ws.Cells.Item(951, 12).Value = 3500
$ws.Cells.Item(951, 13).Value = 3367
$ws.Cells.Item(951, 16).Value = 337

# Row 952
$ws.Cells.Item(952, 4).Value = 44377
$ws.Cells.Item(952, 11).Value = 2500
$ws.Cells.Item(952, 12).Value = 3000
$ws.Cells.Item(952, 13).Value = 2667
$ws.Cells.Item(952, 16).Value = 267

# Row 953
$ws.Cells.Item(953, 4).Value = 44377
$ws.Cells.Item(953, 10).Value = 350
$ws.Cells.Item(953, 11).Value = 2000
$ws.Cells.Item(953, 12).Value = 2500
$ws.Cells.Item(953, 13).Value = 2286
$ws.Cells.Item(953, 16).Value = 229

# Row 954
$ws.Cells.Item(954, 4).Value = 44397
$ws.Cells.Item(954, 10).Value = 400
$ws.Cells.Item(954, 11).Value = 4000
$ws.Cells.Item(954, 12).Value = 4500
$ws.Cells.Item(954, 13).Value = 4250
$ws.Cells.Item(954, 16).Value = 425

# Row 955
$ws.Cells.Item(955, 4).Value = 44397
$ws.Cells.Item(955, 10).Value = 450
$ws.Cells.Item(955, 11).Value = 3500
$ws.Cells.Item(955, 12).Value = 4000
$ws.Cells.Item(955, 13).Value = 3750
$ws.Cells.Item(955, 16).Value = 375

# Row 956
$ws.Cells.Item(956, 4).Value = 44397
$ws.Cells.Item(956, 10).Value = 500
$ws.Cells.Item(956, 11).Value = 3000
$ws.Cells.Item(956, 12).Value = 3500
$ws.Cells.Item(956, 13).Value = 3250
$ws.Cells.Item(956, 16).Value = 325

# Row 957
$ws.Cells.Item(957, 4).Value = 44181
$ws.Cells.Item(957, 10).Value = 360
$ws.Cells.Item(957, 11).Value = 2500
$ws.Cells.Item(957, 12).Value = 3000
$ws.Cells.Item(957, 13).Value = 2750
$ws.Cells.Item(957, 16).Value = 275

# Row 958
$ws.Cells.Item(958, 4).Value = 44181
$ws.Cells.Item(958, 10).Value = 400
$ws.Cells.Item(958, 11).Value = 2000
$ws.Cells.Item(958, 12).Value = 2500
$ws.Cells.Item(958, 13).Value = 2250
$ws.Cells.Item(958, 16).Value = 225

# Row 959
$ws.Cells.Item(959, 4).Value = 44181
$ws.Cells.Item(959, 10).Value = 450
$ws.Cells.Item(959, 11).Value = 1800
$ws.Cells.Item(959, 12).Value = 2000
$ws.Cells.Item(959, 13).Value = 1900
$ws.Cells.Item(959, 16).Value = 190

# Row 960
$ws.Cells.Item(960, 4).Value = 44497
$ws.Cells.Item(960, 10).Value = 200
$ws.Cells.Item(960, 11).Value = 6000
$ws.Cells.Item(960, 12).Value = 6500
$ws.Cells.Item(960, 13).Value = 6250
$ws.Cells.Item(960, 16).Value = 625

# Row 961
$ws.Cells.Item(961, 4).Value = 44497
$ws.Cells.Item(961, 10).Value = 250
$ws.Cells.Item(961, 11).Value = 5000
$ws.Cells.Item(961, 12).Value = 5500
$ws.Cells.Item(961, 13).Value = 5250
$ws.Cells.Item(961, 16).Value = 525

# Row 962
$ws.Cells.Item(962, 4).Value = 44497
$ws.Cells.Item(962, 11).Value = 4000
$ws.Cells.Item(962, 12).Value = 4500
$ws.Cells.Item(962, 13).Value = 4250
$ws.Cells.Item(962, 16).Value = 425

# Row 963
$ws.Cells.Item(963, 4).Value = 44285
$ws.Cells.Item(963, 10).Value = 300
$ws.Cells.Item(963, 11).Value = 4000
$ws.Cells.Item(963, 12).Value = 4500
$ws.Cells.Item(963, 13).Value = 4250
$ws.Cells.Item(963, 16).Value = 425

# Row 964
$ws.Cells.Item(964, 4).Value = 44285
$ws.Cells.Item(964, 10).Value = 360
$ws.Cells.Item(964, 11).Value = 3000
$ws.Cells.Item(964, 12).Value = 3500
$ws.Cells.Item(964, 13).Value = 3250
$ws.Cells.Item(964, 16).Value = 325

# Row 965
$ws.Cells.Item(965, 4).Value = 44285
$ws.Cells.Item(965, 11).Value = 2500
$ws.Cells.Item(965, 12).Value = 3000
$ws.Cells.Item(965, 13).Value = 2750
$ws.Cells.Item(965, 16).Value = 275

# Row 966
$ws.Cells.Item(966, 4).Value = 44362
$ws.Cells.Item(966, 10).Value = 270
$ws.Cells.Item(966, 11).Value = 4500
$ws.Cells.Item(966, 12).Value = 5000
$ws.Cells.Item(966, 13).Value = 4750
$ws.Cells.Item(966, 16).Value = 475

# Row 967
$ws.Cells.Item(967, 4).Value = 44362
$ws.Cells.Item(967, 11).Value = 4000
$ws.Cells.Item(967, 12).Value = 4500
$ws.Cells.Item(967, 13).Value = 4250
$ws.Cells.Item(967, 16).Value = 425

# Row 968
$ws.Cells.Item(968, 4).Value = 44362
$ws.Cells.Item(968, 11).Value = 3500
$ws.Cells.Item(968, 12).Value = 4000
$ws.Cells.Item(968, 13).Value = 3750
$ws.Cells.Item(968, 16).Value = 375

# Row 969
$ws.Cells.Item(969, 4).Value = 44557
$ws.Cells.Item(969, 10).Value = 250
$ws.Cells.Item(969, 11).Value = 4000
$ws.Cells.Item(969, 12).Value = 4500
$ws.Cells.Item(969, 13).Value = 4250
$ws.Cells.Item(969, 16).Value = 425

# Row 970
$ws.Cells.Item(970, 4).Value = 44557
$ws.Cells.Item(970, 10).Value = 300
$ws.Cells.Item(970, 11).Value = 3000
$ws.Cells.Item(970, 12).Value = 3500
$ws.Cells.Item(970, 13).Value = 3250
$ws.Cells.Item(970, 16).Value = 325

# Row 971
$ws.Cells.Item(971, 4).Value = 44557
$ws.Cells.Item(971, 10).Value = 300
$ws.Cells.Item(971, 11).Value = 2500
$ws.Cells.Item(971, 12).Value = 3000
$ws.Cells.Item(971, 13).Value = 2750
$ws.Cells.Item(971, 16).Value = 275

# Row 972
$ws.Cells.Item(972, 4).Value = 44747
$ws.Cells.Item(972, 10).Value = 300
$ws.Cells.Item(972, 11).Value = 3000
$ws.Cells.Item(972, 12).Value = 3500
$ws.Cells.Item(972, 13).Value = 3250
$ws.Cells.Item(972, 16).Value = 325

# Row 973
$ws.Cells.Item(973, 4).Value = 44747
$ws.Cells.Item(973, 10).Value = 350
$ws.Cells.Item(973, 11).Value = 2500
$ws.Cells.Item(973, 12).Value = 3000
$ws.Cells.Item(973, 13).Value = 2750
$ws.Cells.Item(973, 16).Value = 275

# Row 974
$ws.Cells.Item(974, 4).Value = 44747
$ws.Cells.Item(974, 10).Value = 400

# Row 975
$ws.Cells.Item(975, 4).Value = 44357
$ws.Cells.Item(975, 10).Value = 340

# Row 976
$ws.Cells.Item(976, 4).Value = 44357

# Row 977
$ws.Cells.Item(977, 4).Value = 44357
$ws.Cells.Item(977, 10).Value = 500
$ws.Cells.Item(977, 11).Value = 2000
$ws.Cells.Item(977, 12).Value = 2500
$ws.Cells.Item(977, 13).Value = 2250
$ws.Cells.Item(977, 16).Value = 225

# Row 978
$ws.Cells.Item(978, 4).Value = 44279
$ws.Cells.Item(978, 10).Value = 300
$ws.Cells.Item(978, 11).Value = 3500
$ws.Cells.Item(978, 12).Value = 4000
$ws.Cells.Item(978, 13).Value = 3750
$ws.Cells.Item(978, 16).Value = 375

# Row 979
$ws.Cells.Item(979, 4).Value = 44279
$ws.Cells.Item(979, 10).Value = 400
$ws.Cells.Item(979, 11).Value = 3000
$ws.Cells.Item(979, 12).Value = 3500
$ws.Cells.Item(979, 13).Value = 3250
$ws.Cells.Item(979, 16).Value = 325

# Row 980
$ws.Cells.Item(980, 4).Value = 44279
$ws.Cells.Item(980, 10).Value = 560
$ws.Cells.Item(980, 11).Value = 2500
$ws.Cells.Item(980, 12).Value = 3000
$ws.Cells.Item(980, 13).Value = 2750
$ws.Cells.Item(980, 16).Value = 275

# Row 981
$ws.Cells.Item(981, 4).Value = 44551
$ws.Cells.Item(981, 11).Value = 3000
$ws.Cells.Item(981, 12).Value = 3500
$ws.Cells.Item(981, 13).Value = 3250
$ws.Cells.Item(981, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(981, 16).Value = 325
$ws.Cells.Item(981, 17).Value = 10

# Row 982
$ws.Cells.Item(982, 4).Value = 44551
$ws.Cells.Item(982, 9).Value = 'Segunda'
$ws.Cells.Item(982, 10).Value = 270
$ws.Cells.Item(982, 11).Value = 2500
$ws.Cells.Item(982, 12).Value = 3000
$ws.Cells.Item(982, 13).Value = 2750
$ws.Cells.Item(982, 16).Value = 275

# Row 983
$ws.Cells.Item(983, 4).Value = 44551
$ws.Cells.Item(983, 9).Value = 'Tercera'
$ws.Cells.Item(983, 11).Value = 2000
$ws.Cells.Item(983, 12).Value = 2500
$ws.Cells.Item(983, 13).Value = 2250
$ws.Cells.Item(983, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(983, 16).Value = 225
$ws.Cells.Item(983, 17).Value = 10

# Row 984
$ws.Cells.Item(984, 9).Value = 'Primera'
$ws.Cells.Item(984, 10).Value = 250
$ws.Cells.Item(984, 11).Value = 5000
$ws.Cells.Item(984, 12).Value = 5500
$ws.Cells.Item(984, 13).Value = 5250
$ws.Cells.Item(984, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(984, 16).Value = 292
$ws.Cells.Item(984, 17).Value = 18

# Row 985
$ws.Cells.Item(985, 9).Value = 'Primera'
$ws.Cells.Item(985, 10).Value = 350
$ws.Cells.Item(985, 11).Value = 2000
$ws.Cells.Item(985, 12).Value = 2500
$ws.Cells.Item(985, 13).Value = 2250
$ws.Cells.Item(985, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(985, 16).Value = 225
$ws.Cells.Item(985, 17).Value = 10

# Row 986
$ws.Cells.Item(986, 9).Value = 'Segunda'
$ws.Cells.Item(986, 10).Value = 300
$ws.Cells.Item(986, 11).Value = 4000
$ws.Cells.Item(986, 12).Value = 4500
$ws.Cells.Item(986, 13).Value = 4250
$ws.Cells.Item(986, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(986, 16).Value = 236
$ws.Cells.Item(986, 17).Value = 18

# Row 987
$ws.Cells.Item(987, 4).Value = 44757
$ws.Cells.Item(987, 9).Value = 'Segunda'
$ws.Cells.Item(987, 10).Value = 450
$ws.Cells.Item(987, 11).Value = 1500
$ws.Cells.Item(987, 12).Value = 2000
$ws.Cells.Item(987, 13).Value = 1750
$ws.Cells.Item(987, 16).Value = 175

# Row 988
$ws.Cells.Item(988, 4).Value = 44757
$ws.Cells.Item(988, 9).Value = 'Tercera'
$ws.Cells.Item(988, 10).Value = 300
$ws.Cells.Item(988, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(988, 16).Value = 181
$ws.Cells.Item(988, 17).Value = 18

# Row 989
$ws.Cells.Item(989, 4).Value = 44757
$ws.Cells.Item(989, 10).Value = 450
$ws.Cells.Item(989, 11).Value = 1000
$ws.Cells.Item(989, 12).Value = 1500
$ws.Cells.Item(989, 13).Value = 1250
$ws.Cells.Item(989, 16).Value = 125

# Row 990
$ws.Cells.Item(990, 4).Value = 44321
$ws.Cells.Item(990, 10).Value = 400
$ws.Cells.Item(990, 11).Value = 3500
$ws.Cells.Item(990, 12).Value = 4000
$ws.Cells.Item(990, 13).Value = 3750
$ws.Cells.Item(990, 16).Value = 375

# Row 991
$ws.Cells.Item(991, 4).Value = 44321
$ws.Cells.Item(991, 10).Value = 450
$ws.Cells.Item(991, 11).Value = 3000
$ws.Cells.Item(991, 12).Value = 3500
$ws.Cells.Item(991, 13).Value = 3250
$ws.Cells.Item(991, 16).Value = 325

# Row 992
$ws.Cells.Item(992, 4).Value = 44321
$ws.Cells.Item(992, 10).Value = 500
$ws.Cells.Item(992, 11).Value = 2500
$ws.Cells.Item(992, 12).Value = 3000
$ws.Cells.Item(992, 13).Value = 2750
$ws.Cells.Item(992, 16).Value = 275

# Row 993
$ws.Cells.Item(993, 4).Value = 44438
$ws.Cells.Item(993, 10).Value = 300
$ws.Cells.Item(993, 11).Value = 5000
$ws.Cells.Item(993, 12).Value = 5500
$ws.Cells.Item(993, 13).Value = 5250
$ws.Cells.Item(993, 16).Value = 525

# Row 994
$ws.Cells.Item(994, 4).Value = 44438
$ws.Cells.Item(994, 10).Value = 350
$ws.Cells.Item(994, 11).Value = 4500
$ws.Cells.Item(994, 12).Value = 5000
$ws.Cells.Item(994, 13).Value = 4750
$ws.Cells.Item(994, 16).Value = 475

# Row 995
$ws.Cells.Item(995, 4).Value = 44438
$ws.Cells.Item(995, 10).Value = 450
$ws.Cells.Item(995, 11).Value = 3500
$ws.Cells.Item(995, 12).Value = 4000
$ws.Cells.Item(995, 13).Value = 3750
$ws.Cells.Item(995, 16).Value = 375

# Row 996
$ws.Cells.Item(996, 4).Value = 44355
$ws.Cells.Item(996, 10).Value = 450
$ws.Cells.Item(996, 11).Value = 3500
$ws.Cells.Item(996, 12).Value = 4000
$ws.Cells.Item(996, 13).Value = 3750
$ws.Cells.Item(996, 16).Value = 375

# Row 997
$ws.Cells.Item(997, 4).Value = 44355
$ws.Cells.Item(997, 10).Value = 500
$ws.Cells.Item(997, 11).Value = 3000
$ws.Cells.Item(997, 12).Value = 3500
$ws.Cells.Item(997, 13).Value = 3250
$ws.Cells.Item(997, 16).Value = 325

# Row 998
$ws.Cells.Item(998, 4).Value = 44355
$ws.Cells.Item(998, 10).Value = 560
$ws.Cells.Item(998, 11).Value = 2500
$ws.Cells.Item(998, 12).Value = 3000
$ws.Cells.Item(998, 13).Value = 2750
$ws.Cells.Item(998, 16).Value = 275

# Row 999
$ws.Cells.Item(999, 4).Value = 44657
$ws.Cells.Item(999, 10).Value = 300
$ws.Cells.Item(999, 11).Value = 8000
$ws.Cells.Item(999, 12).Value = 8500
$ws.Cells.Item(999, 13).Value = 8250
$ws.Cells.Item(999, 16).Value = 825

# Row 1000
$ws.Cells.Item(1000, 4).Value = 44657
$ws.Cells.Item(1000, 10).Value = 355
$ws.Cells.Item(1000, 11).Value = 7000
$ws.Cells.Item(1000, 12).Value = 7500
$ws.Cells.Item(1000, 13).Value = 7254
$ws.Cells.Item(1000, 16).Value = 725

# Row 1001
$ws.Cells.Item(1001, 4).Value = 44657
$ws.Cells.Item(1001, 11).Value = 5500
$ws.Cells.Item(1001, 12).Value = 6000
$ws.Cells.Item(1001, 13).Value = 5750
$ws.Cells.Item(1001, 16).Value = 575

# Row 1002
$ws.Cells.Item(1002, 4).Value = 44391
$ws.Cells.Item(1002, 10).Value = 250
$ws.Cells.Item(1002, 11).Value = 4000
$ws.Cells.Item(1002, 12).Value = 4500
$ws.Cells.Item(1002, 13).Value = 4250
$ws.Cells.Item(1002, 16).Value = 425

# Row 1003
$ws.Cells.Item(1003, 4).Value = 44391
$ws.Cells.Item(1003, 10).Value = 400
$ws.Cells.Item(1003, 11).Value = 3500
$ws.Cells.Item(1003, 12).Value = 4000
$ws.Cells.Item(1003, 13).Value = 3750
$ws.Cells.Item(1003, 16).Value = 375

# Row 1004
$ws.Cells.Item(1004, 4).Value = 44391
$ws.Cells.Item(1004, 11).Value = 3000
$ws.Cells.Item(1004, 12).Value = 3500
$ws.Cells.Item(1004, 13).Value = 3250
$ws.Cells.Item(1004, 16).Value = 325

# Row 1005
$ws.Cells.Item(1005, 4).Value = 44453
$ws.Cells.Item(1005, 10).Value = 300
$ws.Cells.Item(1005, 11).Value = 6000
$ws.Cells.Item(1005, 12).Value = 6500
$ws.Cells.Item(1005, 13).Value = 6250
$ws.Cells.Item(1005, 16).Value = 625

# Row 1006
$ws.Cells.Item(1006, 4).Value = 44453
$ws.Cells.Item(1006, 10).Value = 350
$ws.Cells.Item(1006, 11).Value = 5000
$ws.Cells.Item(1006, 12).Value = 5500
$ws.Cells.Item(1006, 13).Value = 5250
$ws.Cells.Item(1006, 16).Value = 525

# Row 1007
$ws.Cells.Item(1007, 4).Value = 44453
$ws.Cells.Item(1007, 10).Value = 450
$ws.Cells.Item(1007, 11).Value = 4000
$ws.Cells.Item(1007, 12).Value = 4500
$ws.Cells.Item(1007, 13).Value = 4250
$ws.Cells.Item(1007, 16).Value = 425

# Row 1008
$ws.Cells.Item(1008, 4).Value = 44609
$ws.Cells.Item(1008, 10).Value = 250
$ws.Cells.Item(1008, 11).Value = 4000
$ws.Cells.Item(1008, 12).Value = 4500
$ws.Cells.Item(1008, 13).Value = 4250
$ws.Cells.Item(1008, 16).Value = 425

# Row 1009
$ws.Cells.Item(1009, 4).Value = 44609
$ws.Cells.Item(1009, 10).Value = 300
$ws.Cells.Item(1009, 11).Value = 3500
$ws.Cells.Item(1009, 13).Value = 3750
$ws.Cells.Item(1009, 16).Value = 375

# Row 1010
$ws.Cells.Item(1010, 4).Value = 44609
$ws.Cells.Item(1010, 10).Value = 270

# Row 1011
$ws.Cells.Item(1011, 4).Value = 44358
$ws.Cells.Item(1011, 10).Value = 340
$ws.Cells.Item(1011, 11).Value = 4500
$ws.Cells.Item(1011, 12).Value = 5000
$ws.Cells.Item(1011, 13).Value = 4750
$ws.Cells.Item(1011, 16).Value = 475

# Row 1012
$ws.Cells.Item(1012, 4).Value = 44358
$ws.Cells.Item(1012, 10).Value = 225
$ws.Cells.Item(1012, 11).Value = 4000
$ws.Cells.Item(1012, 12).Value = 4000
$ws.Cells.Item(1012, 13).Value = 4000
$ws.Cells.Item(1012, 16).Value = 400

# Row 1013
$ws.Cells.Item(1013, 4).Value = 44358
$ws.Cells.Item(1013, 10).Value = 450
$ws.Cells.Item(1013, 11).Value = 3000
$ws.Cells.Item(1013, 12).Value = 3500
$ws.Cells.Item(1013, 13).Value = 3250
$ws.Cells.Item(1013, 16).Value = 325

# --- Append new rows 1014-1016 (continuation of the shifted block) ---
# Row 1014
$ws.Cells.Item(1014, 1).Value = 1
$ws.Cells.Item(1014, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(1014, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(1014, 4).Value = 44572
$ws.Cells.Item(1014, 5).Value = 15
$ws.Cells.Item(1014, 6).Value = 100112020
$ws.Cells.Item(1014, 7).Value = 'Tomate'
$ws.Cells.Item(1014, 8).Value = 'Larga vida'
$ws.Cells.Item(1014, 9).Value = 'Primera'
$ws.Cells.Item(1014, 10).Value = 270
$ws.Cells.Item(1014, 11).Value = 4000
$ws.Cells.Item(1014, 12).Value = 4500
$ws.Cells.Item(1014, 13).Value = 4250
$ws.Cells.Item(1014, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(1014, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(1014, 16).Value = 425
$ws.Cells.Item(1014, 17).Value = 10
$ws.Cells.Item(1014, 18).Value = 'Hortaliza'
$ws.Cells.Item(1014, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1015
$ws.Cells.Item(1015, 1).Value = 1
$ws.Cells.Item(1015, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(1015, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(1015, 4).Value = 44572
$ws.Cells.Item(1015, 5).Value = 15
$ws.Cells.Item(1015, 6).Value = 100112020
$ws.Cells.Item(1015, 7).Value = 'Tomate'
$ws.Cells.Item(1015, 8).Value = 'Larga vida'
$ws.Cells.Item(1015, 9).Value = 'Segunda'
$ws.Cells.Item(1015, 10).Value = 300
$ws.Cells.Item(1015, 11).Value = 3000
$ws.Cells.Item(1015, 12).Value = 3500
$ws.Cells.Item(1015, 13).Value = 3250
$ws.Cells.Item(1015, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(1015, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(1015, 16).Value = 325
$ws.Cells.Item(1015, 17).Value = 10
$ws.Cells.Item(1015, 18).Value = 'Hortaliza'
$ws.Cells.Item(1015, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1016
$ws.Cells.Item(1016, 1).Value = 1
$ws.Cells.Item(1016, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(1016, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(1016, 4).Value = 44572
$ws.Cells.Item(1016, 5).Value = 15
$ws.Cells.Item(1016, 6).Value = 100112020
$ws.Cells.Item(1016, 7).Value = 'Tomate'
$ws.Cells.Item(1016, 8).Value = 'Larga vida'
$ws.Cells.Item(1016, 9).Value = 'Tercera'
$ws.Cells.Item(1016, 10).Value = 300
$ws.Cells.Item(1016, 11).Value = 2500
$ws.Cells.Item(1016, 12).Value = 3000
$ws.Cells.Item(1016, 13).Value = 2750
$ws.Cells.Item(1016, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(1016, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(1016, 16).Value = 275
$ws.Cells.Item(1016, 17).Value = 10
$ws.Cells.Item(1016, 18).Value = 'Hortaliza'
$ws.Cells.Item(1016, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
